$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.877.55'
$ws.Range("D3").Value = '3.104.24'
$ws.Range("E3").Value = '  +5.28%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '581.68'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.87%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '173.76'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +7.72%  '
$ws.Range("D8").Value = '3.099.33'
$ws.Range("E8").Value = '  +5.19%  '
$ws.Range("E9").Value = '  +1.41%  '
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.156'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.09%  '
$ws.Range("B11").Value = 'Toncoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.44'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.53%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.485'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.21%  '
$ws.Range("E13").Value = '  +2.31%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.48'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +7.88%  '
$ws.Range("D16").Value = '3.615.04'
$ws.Range("E16").Value = '  +5.08%  '
$ws.Range("D17").Value = '66.882.70'
$ws.Range("E17").Value = '  +2.00%  '
$ws.Range("E18").Value = '  +1.72%  '
$ws.Range("D19").Value = '3.108.87'
$ws.Range("E19").Value = '  +5.37%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.14'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.91%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '480.75'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +7.67%  '
$ws.Range("E22").Value = '  +2.94%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.53'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.29%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.04'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.41%  '
$ws.Range("E25").Value = '  +6.57%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.15'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +7.45%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.07'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.59%  '
$ws.Range("E28").Value = '  +0.03%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.36%  '
$ws.Range("E30").Value = '  -4.16%  '
$ws.Range("E31").Value = '  +3.82%  '
$ws.Range("E32").Value = '  -0.53%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '28.70'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.68%  '
$ws.Range("E34").Value = '  +1.03%  '
$ws.Range("E35").Value = '  +0.08%  '
$ws.Range("E36").Value = '  +3.57%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.994'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.36%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '48.13'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.82%  '
$ws.Range("E39").Value = '  +7.49%  '
$ws.Range("E40").Value = '  +5.14%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '50.13'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.10%  '
$ws.Range("E42").Value = '  +0.47%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.67'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.93%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.83'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.57%  '
$ws.Range("E45").Value = '  +3.23%  '
$ws.Range("D46").Value = '2.827.25'
$ws.Range("E46").Value = '  +5.52%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '384.03'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.60%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '134.63'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.12%  '
$ws.Range("E49").Value = '  -0.01%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '24.90'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.64%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.24'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.30%  '
